$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain plain text so values like "1.00" or
# "3.20" keep their exact formatting instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '51.569.31'
$ws.Range('E2').Value = '  -0.36%  '

$ws.Range('D3').Value = '2.939.21'
$ws.Range('E3').Value = '  +0.86%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.21%  '

$ws.Range('D5').Value = '358.16'
$ws.Range('E5').Value = '  +0.65%  '

$ws.Range('E6').Value = '  -3.77%  '

$ws.Range('E7').Value = '  -2.91%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  -5.20%  '

$ws.Range('D10').Value = '37.28'
$ws.Range('E10').Value = '  -4.47%  '

$ws.Range('E11').Value = '  +2.25%  '

$ws.Range('D12').Value = '0.0844'
$ws.Range('E12').Value = '  -2.96%  '

$ws.Range('D13').Value = '18.74'
$ws.Range('E13').Value = '  -3.95%  '

$ws.Range('D14').Value = '3.404.16'
$ws.Range('E14').Value = '  +0.77%  '

$ws.Range('E15').Value = '  -4.93%  '

$ws.Range('D16').Value = '2.942.02'
$ws.Range('E16').Value = '  +1.29%  '

$ws.Range('D17').Value = '0.975'
$ws.Range('E17').Value = '  -0.59%  '

$ws.Range('D18').Value = '51.488.52'
$ws.Range('E18').Value = '  -0.61%  '

$ws.Range('E19').Value = '  -1.29%  '

$ws.Range('D20').Value = '7.31'
$ws.Range('E20').Value = '  -2.92%  '

$ws.Range('D21').Value = '13.18'
$ws.Range('E21').Value = '  -4.39%  '

$ws.Range('D22').Value = '0.0₃0952'
$ws.Range('E22').Value = '  -2.50%  '

$ws.Range('D23').Value = '68.94'
$ws.Range('E23').Value = '  -2.55%  '

$ws.Range('D24').Value = '262.96'
$ws.Range('E24').Value = '  -2.14%  '

$ws.Range('D25').Value = '2.69'
$ws.Range('E25').Value = '  -4.18%  '

$ws.Range('D26').Value = '0.176'
$ws.Range('E26').Value = '  -5.44%  '

$ws.Range('D27').Value = '26.37'
$ws.Range('E27').Value = '  -1.76%  '

$ws.Range('E28').Value = '  +0.08%  '

$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '0.109'
$ws.Range('E29').Value = '  +1.90%  '

$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '7.16'
$ws.Range('E30').Value = '  -5.43%  '

$ws.Range('D31').Value = '6.26'
$ws.Range('E31').Value = '  +3.57%  '

$ws.Range('D32').Value = '10.08'
$ws.Range('E32').Value = '  -4.08%  '

$ws.Range('E33').Value = '  +4.14%  '

$ws.Range('D34').Value = '35.26'
$ws.Range('E34').Value = '  -6.48%  '

$ws.Range('D35').Value = '51.35'
$ws.Range('E35').Value = '  -1.84%  '

$ws.Range('E36').Value = '  +0.29%  '

$ws.Range('D37').Value = '0.0426'
$ws.Range('E37').Value = '  -3.12%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '2.84'
$ws.Range('E38').Value = '  +4.65%  '

$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '3.20'
$ws.Range('E39').Value = '  +0.07%  '

$ws.Range('D40').Value = '17.18'
$ws.Range('E40').Value = '  -5.73%  '

$ws.Range('E41').Value = '  -4.96%  '

$ws.Range('E42').Value = '  -3.70%  '

$ws.Range('D43').Value = '23.01'
$ws.Range('E43').Value = '  +0.25%  '

$ws.Range('D44').Value = '120.72'
$ws.Range('E44').Value = '  +1.45%  '

$ws.Range('E45').Value = '  -1.82%  '

$ws.Range('D46').Value = '2.086.04'
$ws.Range('E46').Value = '  -1.56%  '

$ws.Range('E47').Value = '  -6.58%  '

$ws.Range('E48').Value = '  -7.19%  '

$ws.Range('D49').Value = '3.229.75'
$ws.Range('E49').Value = '  +0.69%  '

$ws.Range('D50').Value = '0.237'
$ws.Range('E50').Value = '  -4.98%  '

$ws.Range('D51').Value = '0.0316'
$ws.Range('E51').Value = '  -4.94%  '
